# Apply the "pre-week data" update to the sample-size workbook.
# Sheet "Level0" (sheet1): a handful of cell values change (redistribution
#   between the CTRL/MM sub-columns and their row/column totals) while the
#   category rows themselves stay the same.
# Sheet "Level1" (sheet2): the cell-type breakdown is refreshed -- two new
#   categories appear (EN_Lymph, GC_Atretic), one disappears (GC_Mural), and
#   every remaining value is updated to the new counts.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Level0")
$ws2 = $wb.Worksheets.Item("Level1")

# --- Sheet1 (Level0): update changed cells ---
$ws1.Range("G2").Value = 28
$ws1.Range("I2").Value = 49
$ws1.Range("J2").Value = 83
$ws1.Range("B3").Value = 9
$ws1.Range("C3").Value = 9
$ws1.Range("D3").Value = 14
$ws1.Range("F3").Value = 11
$ws1.Range("H3").Value = 32
$ws1.Range("I3").Value = 31
$ws1.Range("J3").Value = 63
$ws1.Range("B4").Value = 864
$ws1.Range("C4").Value = 1069
$ws1.Range("D4").Value = 845
$ws1.Range("E4").Value = 839
$ws1.Range("F4").Value = 1026
$ws1.Range("G4").Value = 857
$ws1.Range("H4").Value = 2778
$ws1.Range("I4").Value = 2722
$ws1.Range("J4").Value = 5500
$ws1.Range("B6").Value = 727
$ws1.Range("C6").Value = 512
$ws1.Range("D6").Value = 639
$ws1.Range("E6").Value = 535
$ws1.Range("F6").Value = 717
$ws1.Range("G6").Value = 505
$ws1.Range("H6").Value = 1878
$ws1.Range("I6").Value = 1757
$ws1.Range("J6").Value = 3635

# --- Sheet2 (Level1): structural edits ---
$ws2.Rows.Item(3).Insert()   # insert row for EN_Lymph
$ws2.Rows.Item(7).Insert()   # insert row for GC_Atretic
$ws2.Rows.Item(10).Delete()  # remove GC_Mural row

# --- Sheet2 (Level1): write full updated table ---
$ws2.Range("A2").Value = "EN_Blood"
$ws2.Range("B2").Value = 13
$ws2.Range("C2").Value = 8
$ws2.Range("D2").Value = 13
$ws2.Range("E2").Value = 13
$ws2.Range("F2").Value = 7
$ws2.Range("G2").Value = 28
$ws2.Range("H2").Value = 34
$ws2.Range("I2").Value = 48
$ws2.Range("J2").Value = 82
$ws2.Range("A3").Value = "EN_Lymph"
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 1
$ws2.Range("J3").Value = 1
$ws2.Range("A4").Value = "Epithelium"
$ws2.Range("B4").Value = 9
$ws2.Range("C4").Value = 10
$ws2.Range("D4").Value = 14
$ws2.Range("E4").Value = 10
$ws2.Range("F4").Value = 11
$ws2.Range("G4").Value = 10
$ws2.Range("H4").Value = 33
$ws2.Range("I4").Value = 31
$ws2.Range("J4").Value = 64
$ws2.Range("A5").Value = "GC_Active CL"
$ws2.Range("B5").Value = 9
$ws2.Range("C5").Value = 3
$ws2.Range("D5").Value = 8
$ws2.Range("E5").Value = 2
$ws2.Range("F5").Value = 3
$ws2.Range("G5").Value = 9
$ws2.Range("H5").Value = 20
$ws2.Range("I5").Value = 14
$ws2.Range("J5").Value = 34
$ws2.Range("A6").Value = "GC_Antral"
$ws2.Range("B6").Value = 459
$ws2.Range("C6").Value = 680
$ws2.Range("D6").Value = 532
$ws2.Range("E6").Value = 450
$ws2.Range("F6").Value = 590
$ws2.Range("G6").Value = 528
$ws2.Range("H6").Value = 1671
$ws2.Range("I6").Value = 1568
$ws2.Range("J6").Value = 3239
$ws2.Range("A7").Value = "GC_Atretic"
$ws2.Range("B7").Value = 27
$ws2.Range("C7").Value = 18
$ws2.Range("D7").Value = 15
$ws2.Range("E7").Value = 20
$ws2.Range("F7").Value = 20
$ws2.Range("G7").Value = 21
$ws2.Range("H7").Value = 60
$ws2.Range("I7").Value = 61
$ws2.Range("J7").Value = 121
$ws2.Range("A8").Value = "GC_Luteinizing"
$ws2.Range("B8").Value = 186
$ws2.Range("C8").Value = 110
$ws2.Range("D8").Value = 102
$ws2.Range("E8").Value = 150
$ws2.Range("F8").Value = 84
$ws2.Range("G8").Value = 74
$ws2.Range("H8").Value = 398
$ws2.Range("I8").Value = 308
$ws2.Range("J8").Value = 706
$ws2.Range("A9").Value = "GC_Mitotic"
$ws2.Range("B9").Value = 112
$ws2.Range("C9").Value = 208
$ws2.Range("D9").Value = 124
$ws2.Range("E9").Value = 168
$ws2.Range("F9").Value = 259
$ws2.Range("G9").Value = 181
$ws2.Range("H9").Value = 444
$ws2.Range("I9").Value = 608
$ws2.Range("J9").Value = 1052
$ws2.Range("A10").Value = "GC_Preantral"
$ws2.Range("B10").Value = 52
$ws2.Range("C10").Value = 32
$ws2.Range("D10").Value = 48
$ws2.Range("E10").Value = 43
$ws2.Range("F10").Value = 60
$ws2.Range("G10").Value = 39
$ws2.Range("H10").Value = 132
$ws2.Range("I10").Value = 142
$ws2.Range("J10").Value = 274
$ws2.Range("A11").Value = "GC_Regressing CL"
$ws2.Range("B11").Value = 9
$ws2.Range("C11").Value = 12
$ws2.Range("D11").Value = 5
$ws2.Range("E11").Value = 6
$ws2.Range("F11").Value = 7
$ws2.Range("G11").Value = 1
$ws2.Range("H11").Value = 26
$ws2.Range("I11").Value = 14
$ws2.Range("J11").Value = 40
$ws2.Range("A12").Value = "I_Macrophage"
$ws2.Range("B12").Value = 15
$ws2.Range("C12").Value = 13
$ws2.Range("D12").Value = 10
$ws2.Range("E12").Value = 12
$ws2.Range("F12").Value = 10
$ws2.Range("G12").Value = 16
$ws2.Range("H12").Value = 38
$ws2.Range("I12").Value = 38
$ws2.Range("J12").Value = 76
$ws2.Range("A13").Value = "M_Early Theca"
$ws2.Range("B13").Value = 113
$ws2.Range("C13").Value = 126
$ws2.Range("D13").Value = 111
$ws2.Range("E13").Value = 102
$ws2.Range("F13").Value = 179
$ws2.Range("G13").Value = 89
$ws2.Range("H13").Value = 350
$ws2.Range("I13").Value = 370
$ws2.Range("J13").Value = 720
$ws2.Range("A14").Value = "M_Fibroblast-like Stroma"
$ws2.Range("B14").Value = 21
$ws2.Range("C14").Value = 12
$ws2.Range("D14").Value = 28
$ws2.Range("E14").Value = 8
$ws2.Range("F14").Value = 19
$ws2.Range("G14").Value = 17
$ws2.Range("H14").Value = 61
$ws2.Range("I14").Value = 44
$ws2.Range("J14").Value = 105
$ws2.Range("A15").Value = "M_Pericyte"
$ws2.Range("B15").Value = 22
$ws2.Range("C15").Value = 26
$ws2.Range("D15").Value = 42
$ws2.Range("E15").Value = 12
$ws2.Range("F15").Value = 21
$ws2.Range("G15").Value = 22
$ws2.Range("H15").Value = 90
$ws2.Range("I15").Value = 55
$ws2.Range("J15").Value = 145
$ws2.Range("A16").Value = "M_Smooth Muscle"
$ws2.Range("B16").Value = 12
$ws2.Range("C16").Value = 5
$ws2.Range("D16").Value = 9
$ws2.Range("E16").Value = 11
$ws2.Range("F16").Value = 13
$ws2.Range("G16").Value = 11
$ws2.Range("H16").Value = 26
$ws2.Range("I16").Value = 35
$ws2.Range("J16").Value = 61
$ws2.Range("A17").Value = "M_Steroidogenic Stroma"
$ws2.Range("B17").Value = 419
$ws2.Range("C17").Value = 249
$ws2.Range("D17").Value = 325
$ws2.Range("E17").Value = 276
$ws2.Range("F17").Value = 359
$ws2.Range("G17").Value = 226
$ws2.Range("H17").Value = 993
$ws2.Range("I17").Value = 861
$ws2.Range("J17").Value = 1854
$ws2.Range("A18").Value = "M_Steroidogenic Theca"
$ws2.Range("B18").Value = 150
$ws2.Range("C18").Value = 99
$ws2.Range("D18").Value = 135
$ws2.Range("E18").Value = 126
$ws2.Range("F18").Value = 129
$ws2.Range("G18").Value = 144
$ws2.Range("H18").Value = 384
$ws2.Range("I18").Value = 399
$ws2.Range("J18").Value = 783
$ws2.Range("A19").Value = "Total"
$ws2.Range("B19").Value = 1628
$ws2.Range("C19").Value = 1611
$ws2.Range("D19").Value = 1521
$ws2.Range("E19").Value = 1410
$ws2.Range("F19").Value = 1771
$ws2.Range("G19").Value = 1416
$ws2.Range("H19").Value = 4760
$ws2.Range("I19").Value = 4597
$ws2.Range("J19").Value = 9357
